$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (e.g. "43.187.21", "0.510", "9.20") that must
# stay Text, same as the workbook's original inline-string cell typing --
# otherwise plain-decimal-looking values like "0.510"/"9.20"/"1.00" would be
# auto-coerced into Number cells (and lose their trailing zero) on assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.187.21"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "2.317.43"
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "303.77"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").Value = "100.27"
$ws.Range("E6").Value = "  -5.05%  "
$ws.Range("D7").Value = "0.510"
$ws.Range("E7").Value = "  -2.06%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").Value = "34.96"
$ws.Range("E10").Value = "  -3.60%  "
$ws.Range("D11").Value = "51.39"
$ws.Range("E11").Value = "  -3.54%  "
$ws.Range("D12").Value = "0.0793"
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "6.78"
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("D15").Value = "2.698.57"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").Value = "15.63"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "2.327.00"
$ws.Range("E17").Value = "  -2.02%  "
$ws.Range("D18").Value = "0.801"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").Value = "43.263.73"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "11.75"
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("D21").Value = "0.0₃0903"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").Value = "6.08"
$ws.Range("E22").Value = "  -3.37%  "
$ws.Range("D23").Value = "67.59"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").Value = "237.31"
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").Value = "1.97"
$ws.Range("E25").Value = "  -4.23%  "
$ws.Range("D26").Value = "2.53"
$ws.Range("E26").Value = "  -3.40%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "24.86"
$ws.Range("E28").Value = "  -3.86%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "34.45"
$ws.Range("E30").Value = "  -6.61%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "164.71"
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "9.20"
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").Value = "5.04"
$ws.Range("E34").Value = "  -4.45%  "
$ws.Range("E35").Value = "  -4.87%  "
$ws.Range("D36").Value = "4.47"
$ws.Range("E36").Value = "  -6.17%  "
$ws.Range("D37").Value = "0.0703"
$ws.Range("E37").Value = "  -5.63%  "
$ws.Range("D38").Value = "16.77"
$ws.Range("E38").Value = "  -8.50%  "
$ws.Range("D39").Value = "2.89"
$ws.Range("E39").Value = "  -7.34%  "
$ws.Range("D40").Value = "1.81"
$ws.Range("E40").Value = "  -6.86%  "
$ws.Range("D41").Value = "0.102"
$ws.Range("E41").Value = "  -3.46%  "
$ws.Range("E42").Value = "  -3.05%  "
$ws.Range("D43").Value = "2.41"
$ws.Range("E43").Value = "  -8.95%  "
$ws.Range("D44").Value = "1.975.27"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("D45").Value = "0.0284"
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("D46").Value = "18.69"
$ws.Range("E46").Value = "  -6.39%  "
$ws.Range("D47").Value = "2.92"
$ws.Range("E47").Value = "  -7.53%  "
$ws.Range("D48").Value = "9.84"
$ws.Range("E48").Value = "  -6.83%  "
$ws.Range("D49").Value = "4.87"
$ws.Range("E49").Value = "  +3.23%  "
$ws.Range("D50").Value = "54.91"
$ws.Range("E50").Value = "  -5.69%  "
$ws.Range("D51").Value = "2.561.47"
$ws.Range("E51").Value = "  +0.15%  "
